$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "Sizes of files produced by compiler" table (columns D=Empty, E=Sieve, F=Fibo)
# Values entered as text (e.g. "33 792") to match the space-separated thousands formatting
# used throughout the sheet, in the order the author originally typed them.
$ws.Range("D3").Value = "33 792"
$ws.Range("E3").Value = "37 376"
$ws.Range("F3").Value = "36 352"

$ws.Range("D6").Value = "90 586"
$ws.Range("F6").Value = "136 832"
$ws.Range("E6").Value = "137 520"

$ws.Range("D7").Value = "90 793"

$ws.Range("E4").Value = "66 560"

$ws.Range("E7").Value = "94 043"

$ws.Range("E5").Value = "149 504"
$ws.Range("F5").Value = "148 480"
$ws.Range("D5").Value = "47 104"

$ws.Range("E9").Value = "35 840"
$ws.Range("D9").Value = "35 328"

$ws.Range("F4").Value = "66 048"
$ws.Range("D4").Value = "52 224"

$ws.Range("F7").Value = "136 832"
$ws.Range("F9").Value = "35 328"

# Row 8 (Java) sizes were recorded as plain numbers, not text
$ws.Range("D8").Value = 257
$ws.Range("E8").Value = 907
$ws.Range("F8").Value = 603
